# "control amb una sola capa VMIX"
# Restructure the VMIX header/points columns to a single ("current section")
# layer instead of per-section P1..P6 breakdown, update the leading player
# to SECTION 1 / ALEJANDRO MO, and swap the PLAYER1 ranking rows 4/5 back so
# VINCENT H leads (row 4) ahead of ALEJANDRO MO (row 5).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VMIX"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("VMIX")

# Drop the old per-section P1..P6 detail columns (FN:GM) entirely, which
# also shrinks the sheet's used range from A1:GM2 down to A1:FM2.
$ws.Range("FN1:GM2").EntireColumn.Delete()

# Rebuild the header row from EW1 through FM1 with the new, consolidated
# "current section" + per-section summary columns.
$headers = @(
    "C_PUNTS_SECCIO",
    "C_PUNTS_P1",
    "C_PUNTS_P2",
    "C_PUNTS_P3",
    "C_PUNTS_P4",
    "C_PUNTS_P5",
    "C_PUNTS_P6",
    "C_BANDERA",
    "C_PAIS",
    "C_PLAYER",
    "hashtag",
    "1_PUNTS_SECCIO",
    "2_PUNTS_SECCIO",
    "3_PUNTS_SECCIO",
    "4_PUNTS_SECCIO",
    "5_PUNTS_SECCIO",
    "6_PUNTS_SECCIO"
)
$startCol = 153   # column EW
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $startCol + $i).Value = $headers[$i]
}

# Data row 2: section now in progress is SECTION 1, leader is ALEJANDRO MO.
$ws.Range("F2").Value = "SECTION 1"
$ws.Range("W2").Value = "ESP"
$ws.Range("AE2").Value = "ALEJANDRO MO"

$ws.Range("DG2").Value = 40
$ws.Range("DH2").Value = 60
$ws.Range("DM2").Value = 40
$ws.Range("DT2").Value = 30
$ws.Range("DZ2").Value = 30

# Current-leader ("C_") points/flag/country/player block plus the
# per-section summary totals.
$ws.Range("EW2").Value = 40
$ws.Range("EX2").Value = 10
$ws.Range("EY2").Value = 10
$ws.Range("EZ2").Value = 10
$ws.Range("FA2").Value = 10
$ws.Range("FB2").Value = "-"
$ws.Range("FC2").Value = "-"
$ws.Range("FD2").Value = "C:\TRIAL_2021\VMIX\MATERIAL\BANDERES\esp.png"
$ws.Range("FE2").Value = "ESP"
$ws.Range("FF2").Value = "ALEJANDRO MO"
$ws.Range("FG2").Value = "#TrialVIC_2021"
$ws.Range("FH2").Value = 40
$ws.Range("FI2").Value = 0
$ws.Range("FJ2").Value = 0
$ws.Range("FK2").Value = 0
$ws.Range("FL2").Value = 0
$ws.Range("FM2").Value = 0

# ---------------------------------------------------------------------
# Sheet "PLAYER1"
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("PLAYER1")

# Swap the ranking rows: VINCENT H moves back up to row 4 (rank 2),
# ALEJANDRO MO moves to row 5 (rank 3).
$ws4.Range("B4").Value = 1
$ws4.Range("C4").Value = 1
$ws4.Range("D4").Value = "VINCENT H"
$ws4.Range("E4").Value = "HER"
$ws4.Range("F4").Value = "FRA"
$ws4.Range("G4").Value = "C:\TRIAL_2021\VMIX\MATERIAL\BANDERES\fra.png"

$ws4.Range("B5").Value = 2
$ws4.Range("C5").Value = 7
$ws4.Range("D5").Value = "ALEJANDRO MO"
$ws4.Range("E5").Value = "MON"
$ws4.Range("F5").Value = "ESP"
$ws4.Range("G5").Value = "C:\TRIAL_2021\VMIX\MATERIAL\BANDERES\esp.png"

$ws4.Range("H5:M5").ClearContents()
$ws4.Range("T5").Value = 0
